$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1662.8572
$ws.Range("I38").Value = 1662.8572
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 4988.571599999999
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -4616.571599999999
$ws.Range("N38").ClearContents()

$ws.Range("H39").Value = 417.4375
$ws.Range("I39").Value = 138.5
$ws.Range("J39").Value = 510.41666
$ws.Range("K39").Value = 415.5
$ws.Range("L39").Value = 1531.24998
$ws.Range("M39").Value = -119.5
$ws.Range("N39").Value = -2123.24998

$ws.Range("H53").Value = 444.69232
$ws.Range("I53").Value = 191.33333
$ws.Range("K53").Value = 191.33333
$ws.Range("M53").Value = 445.66667

$ws.Range("H86").Value = 7335.4546
$ws.Range("I86").Value = 7544.727
$ws.Range("J86").Value = 7126.1816
$ws.Range("K86").Value = 7544.727
$ws.Range("L86").Value = 7126.1816
$ws.Range("M86").Value = -6421.727
$ws.Range("N86").Value = -9372.1816

$ws.Range("H89").Value = 7335.4546
$ws.Range("I89").Value = 7544.727
$ws.Range("J89").Value = 7126.1816
$ws.Range("K89").Value = 37723.635
$ws.Range("L89").Value = 35630.908
$ws.Range("M89").Value = -32107.635
$ws.Range("N89").Value = -46862.908

$ws.Range("H138").Value = 1236871
$ws.Range("I138").Value = 826.0909
$ws.Range("J138").Value = 2086651.8
$ws.Range("K138").Value = 2478.2727
$ws.Range("L138").Value = 6259955.4
$ws.Range("M138").Value = 2661.7273
$ws.Range("N138").Value = -6270235.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1149.3334
$ws.Range("I4").Value = 225
$ws.Range("J4").Value = 2998
$ws.Range("K4").Value = 225
$ws.Range("L4").Value = 2998
$ws.Range("M4").Value = -109
$ws.Range("N4").Value = -3230

$ws.Range("H32").Value = 31253458
$ws.Range("I32").Value = 35716950
$ws.Range("K32").Value = 35716950
$ws.Range("M32").Value = -35716663

$ws.Range("H61").Value = 38544428
$ws.Range("I61").Value = 83335580
$ws.Range("K61").Value = 83335580
$ws.Range("M61").Value = -83335368

$ws.Range("H74").Value = 15638267
$ws.Range("I74").Value = 27779812
$ws.Range("J74").Value = 27709.428
$ws.Range("K74").Value = 27779812
$ws.Range("L74").Value = 27709.428
$ws.Range("M74").Value = -27778938
$ws.Range("N74").Value = -29457.428

$ws.Range("H77").Value = 15638267
$ws.Range("I77").Value = 27779812
$ws.Range("J77").Value = 27709.428
$ws.Range("K77").Value = 138899060
$ws.Range("L77").Value = 138547.14
$ws.Range("M77").Value = -138894692
$ws.Range("N77").Value = -147283.14

$ws.Range("H88").Value = 1439.2106
$ws.Range("I88").Value = 1393.1111
$ws.Range("J88").Value = 1480.7
$ws.Range("K88").Value = 1393.1111
$ws.Range("L88").Value = 1480.7
$ws.Range("M88").Value = -987.1111000000001
$ws.Range("N88").Value = -2292.7

$ws.Range("H91").Value = 1439.2106
$ws.Range("I91").Value = 1393.1111
$ws.Range("J91").Value = 1480.7
$ws.Range("K91").Value = 1393.1111
$ws.Range("L91").Value = 1480.7
$ws.Range("M91").Value = 10.88889999999992
$ws.Range("N91").Value = -4288.7

$ws.Range("H102").Value = 11711.238
$ws.Range("I102").Value = 15064.6
$ws.Range("J102").Value = 3327.8333
$ws.Range("K102").Value = 15064.6
$ws.Range("L102").Value = 3327.8333
$ws.Range("M102").Value = -13442.6
$ws.Range("N102").Value = -6571.8333

$ws.Range("H122").Value = 2143.875
$ws.Range("I122").Value = 1736
$ws.Range("J122").Value = 4999
$ws.Range("K122").Value = 5208
$ws.Range("L122").Value = 14997
$ws.Range("M122").Value = -2758
$ws.Range("N122").Value = -19897

$ws.Range("H134").Value = 220666.67
$ws.Range("J134").Value = 220666.67
$ws.Range("L134").Value = 220666.67
$ws.Range("N134").Value = -230806.67

$ws.Range("H136").Value = 38544428
$ws.Range("I136").Value = 83335580
$ws.Range("K136").Value = 250006740
$ws.Range("M136").Value = -250004190

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 33099.8
$ws.Range("I82").Value = 7749.5
$ws.Range("K82").Value = 7749.5
$ws.Range("M82").Value = -7366.5

$ws.Range("H85").Value = 33099.8
$ws.Range("I85").Value = 7749.5
$ws.Range("K85").Value = 7749.5
$ws.Range("M85").Value = -6423.5

$ws.Range("H94").Value = 1418
$ws.Range("I94").Value = 1517.8667
$ws.Range("J94").Value = 669
$ws.Range("K94").Value = 1517.8667
$ws.Range("L94").Value = 669
$ws.Range("M94").Value = -1066.8667
$ws.Range("N94").Value = -1571

$ws.Range("H99").Value = 2442.3809
$ws.Range("I99").Value = 1980
$ws.Range("J99").Value = 3367.1428
$ws.Range("K99").Value = 1980
$ws.Range("L99").Value = 3367.1428
$ws.Range("M99").Value = -482
$ws.Range("N99").Value = -6363.1428

$ws.Range("H134").Value = 58010.25
$ws.Range("I134").Value = 6480.75
$ws.Range("K134").Value = 19442.25
$ws.Range("M134").Value = -16907.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 972692.9399999999
$ws.Range("J31").Value = 3331750.2
$ws.Range("L31").Value = 3331750.2
$ws.Range("N31").Value = -3332340.2

$ws.Range("H34").Value = 972692.9399999999
$ws.Range("J34").Value = 3331750.2
$ws.Range("L34").Value = 3331750.2
$ws.Range("N34").Value = -3332154.2

$ws.Range("H134").Value = 457440.9
$ws.Range("I134").Value = 589658.5
$ws.Range("K134").Value = 1768975.5
$ws.Range("M134").Value = -1766440.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 20508.25
$ws.Range("J125").Value = 20508.25
$ws.Range("L125").Value = 61524.75
$ws.Range("N125").Value = -71364.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1867.9474
$ws.Range("J22").Value = 1442.5714
$ws.Range("L22").Value = 1442.5714
$ws.Range("N22").Value = -2032.5714

$ws.Range("H27").Value = 1867.9474
$ws.Range("J27").Value = 1442.5714
$ws.Range("L27").Value = 1442.5714
$ws.Range("N27").Value = -1656.5714

$ws.Range("H40").Value = 4148.625
$ws.Range("I40").Value = 3307.182
$ws.Range("K40").Value = 3307.182
$ws.Range("M40").Value = -3171.182

$ws.Range("H55").Value = 55555740
$ws.Range("I55").Value = 62500172
$ws.Range("J55").Value = 264
$ws.Range("K55").Value = 62500172
$ws.Range("L55").Value = 264
$ws.Range("M55").Value = -62499999
$ws.Range("N55").Value = -610

$ws.Range("H68").Value = 4300.3335
$ws.Range("J68").Value = 9249
$ws.Range("L68").Value = 9249
$ws.Range("N68").Value = -10747

$ws.Range("H71").Value = 4300.3335
$ws.Range("J71").Value = 9249
$ws.Range("L71").Value = 46245
$ws.Range("N71").Value = -53733

$ws.Range("H132").Value = 97860.14
$ws.Range("J132").Value = 127618.5
$ws.Range("L132").Value = 382855.5
$ws.Range("N132").Value = -387915.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 5000
$ws.Range("J20").Value = 5000
$ws.Range("L20").Value = 5000
$ws.Range("N20").Value = -5480

$ws.Range("H75").Value = 10756879
$ws.Range("J75").Value = 10756879
$ws.Range("L75").Value = 10756879
$ws.Range("N75").Value = -10758751

$ws.Range("H78").Value = 10756879
$ws.Range("J78").Value = 10756879
$ws.Range("L78").Value = 32270637
$ws.Range("N78").Value = -32279997

$ws.Range("H113").Value = 1155.8823
$ws.Range("J113").Value = 1145
$ws.Range("L113").Value = 3435
$ws.Range("N113").Value = -7775

$ws.Range("H132").Value = 3898.2646
$ws.Range("I132").Value = 3621.158
$ws.Range("K132").Value = 10863.474
$ws.Range("M132").Value = -8333.474

$ws.Range("H136").Value = 1955.1765
$ws.Range("I136").Value = 1941.3846
$ws.Range("K136").Value = 5824.1538
$ws.Range("M136").Value = -3274.1538
